$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column F (the mostly-empty "Note" column), shifting
# everything to the right of it one column to the left.
$ws.Columns.Item(6).Delete()

# Restore the view state (zoom/selection) as recorded after the edit.
$ws.Application.ActiveWindow.Zoom = 125
$ws.Range("E23").Select()
